$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E16:E19) with the new period values
$ws.Range("E16").Value = "1704"
$ws.Range("E17").Value = "1707"
$ws.Range("E18").Value = "1711"
$ws.Range("E19").Value = "1712"

# Update "Valor Mora" column (F16 and F19) to reflect the swapped amounts
$ws.Range("F16").Value = 29520
$ws.Range("F19").Value = 14760
